$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '62.055.31'
Set-TextCell 'E2' '  +2.61%  '
Set-TextCell 'D3' '2.417.45'
Set-TextCell 'E3' '  +3.53%  '
Set-TextCell 'E4' '  +0.29%  '
Set-TextCell 'D5' '553.82'
Set-TextCell 'E5' '  +2.25%  '
Set-TextCell 'D6' '142.97'
Set-TextCell 'E6' '  +5.51%  '
Set-TextCell 'E7' '  +0.26%  '
Set-TextCell 'D8' '0.532'
Set-TextCell 'E8' '  +2.48%  '
Set-TextCell 'D9' '2.416.89'
Set-TextCell 'E9' '  +3.57%  '
Set-TextCell 'E10' '  +5.33%  '
Set-TextCell 'D12' '5.39'
Set-TextCell 'E12' '  +2.47%  '
Set-TextCell 'D13' '0.352'
Set-TextCell 'E13' '  +4.61%  '
Set-TextCell 'D14' '26.21'
Set-TextCell 'E14' '  +7.66%  '
Set-TextCell 'E15' '  +10.14%  '
Set-TextCell 'D16' '2.857.22'
Set-TextCell 'E16' '  +3.95%  '
Set-TextCell 'D17' '62.021.27'
Set-TextCell 'E17' '  +3.09%  '
Set-TextCell 'D18' '2.416.48'
Set-TextCell 'E18' '  +3.74%  '
Set-TextCell 'D19' '11.11'
Set-TextCell 'E19' '  +5.58%  '
Set-TextCell 'D20' '324.43'
Set-TextCell 'E20' '  +2.61%  '
Set-TextCell 'E21' '  +3.33%  '
Set-TextCell 'D22' '6.74'
Set-TextCell 'E22' '  +3.79%  '
Set-TextCell 'E23' '  +0.24%  '
Set-TextCell 'B24' 'SuiNetwork'
Set-TextCell 'C24' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell 'D24' '1.76'
Set-TextCell 'E24' '  +6.27%  '
Set-TextCell 'B25' 'Litecoin'
Set-TextCell 'C25' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D25' '64.87'
Set-TextCell 'E25' '  +3.57%  '
Set-TextCell 'D26' '9.24'
Set-TextCell 'E26' '  +11.43%  '
Set-TextCell 'D27' '565.24'
Set-TextCell 'E27' '  +15.64%  '
Set-TextCell 'D28' '2.541.86'
Set-TextCell 'D29' '1.00'
Set-TextCell 'E29' '  +0.19%  '
Set-TextCell 'B30' 'InternetComputer(DFINITY)'
Set-TextCell 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D30' '8.37'
Set-TextCell 'E30' '  +6.69%  '
Set-TextCell 'B31' 'PEPE'
Set-TextCell 'C31' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 'D31' '0.0₃0937'
Set-TextCell 'E31' '  +10.59%  '
Set-TextCell 'E32' '  +7.04%  '
Set-TextCell 'D33' '0.148'
Set-TextCell 'E33' '  +3.68%  '
Set-TextCell 'D34' '1.86'
Set-TextCell 'E34' '  +4.68%  '
Set-TextCell 'E35' '  +4.64%  '
Set-TextCell 'D36' '5.73'
Set-TextCell 'E36' '  +11.76%  '
Set-TextCell 'E37' '  +12.30%  '
Set-TextCell 'D38' '1.00'
Set-TextCell 'E38' '  +0.48%  '
Set-TextCell 'E39' '  +6.78%  '
Set-TextCell 'E40' '  +3.46%  '
Set-TextCell 'D41' '18.78'
Set-TextCell 'E41' '  +2.03%  '
Set-TextCell 'D42' '147.46'
Set-TextCell 'E42' '  +4.38%  '
Set-TextCell 'E43' '  +0.31%  '
Set-TextCell 'E44' '  +14.82%  '
Set-TextCell 'E45' '  +8.08%  '
Set-TextCell 'D46' '3.63'
Set-TextCell 'E46' '  +3.24%  '
Set-TextCell 'D47' '0.0539'
Set-TextCell 'E47' '  +6.29%  '
Set-TextCell 'D48' '20.37'
Set-TextCell 'E48' '  +7.91%  '
Set-TextCell 'D49' '0.589'
Set-TextCell 'E49' '  +4.55%  '
Set-TextCell 'D50' '0.0909'
Set-TextCell 'E50' '  +1.84%  '
Set-TextCell 'E51' '  +4.30%  '
